$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column O into the new column P for every row that has
# data in column O, then populate the new 2022 values.
$rows = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17)
foreach ($r in $rows) {
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 11.4
$ws.Range("P6").Value = 12.6
$ws.Range("P7").Value = 9.8000000000000007
$ws.Range("P8").Value = 11.4
$ws.Range("P9").Value = 5.4
$ws.Range("P10").Value = 4.7
$ws.Range("P11").Value = 3.4
$ws.Range("P12").Value = 17.7
$ws.Range("P13").Value = 20.5
$ws.Range("P14").Value = 8.4
$ws.Range("P16").Value = 12.9
$ws.Range("P17").Value = 10.5

# Update the saved selection to match the post-edit workbook state.
$ws.Range("Q4").Select()
